# Updated capital structure database
#
# The underlying dataset changed shape:
#   - Row 2 (previously a standalone/aggregate KASE entry) is refreshed with
#     new historical growth, margin, capital-structure and debt figures.
#   - Row 3 (Kazakhstan Stock Exchange Joint-Stock Company) is refreshed with
#     its own updated figures; its stale "buybacks_cash_returned" cell (T3)
#     no longer applies and is cleared.
#   - Row 4 is a brand new entry for Kaspi.kz Joint Stock Company (LSE:KSPI).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row2Data = @{
    "A2" = 'Kazakhstan'
    "C2" = 'Financial Svcs. (Non-bank & Insurance)'
    "D2" = 0.3
    "E2" = 0.419
    "F2" = 0.28
    "G2" = 0.006390984811366977
    "H2" = 0.006330230279274865
    "I2" = 0.006114649681528662
    "J2" = 0.005314910291875514
    "K2" = 570.63
    "L2" = 0.5591670749632534
    "M2" = 491.9
    "N2" = 0.03817291500143565
    "O2" = 0.8620296864868654
    "P2" = 491.9
    "Q2" = 0.03817291500143565
    "R2" = 0.8620296864868654
    "S2" = 0
    "T2" = 0
    "U2" = 597.3
    "V2" = 0.04635227105175344
    "X2" = 0.02072269771519328
    "Z2" = 1.502502944640754
    "AB2" = 0.02126038812292253
    "AD2" = 489.7
    "AE2" = 0
    "AF2" = 489.7
    "AG2" = -107.6
    "AH2" = 0.03661089430164925
    "AI2" = 0.4115471888393982
    "AJ2" = -0.008420393629925263
    "AK2" = -0.1815727303408707
    "AL2" = 0.082
    "AM2" = -0.174
    "AN2" = 75.80495356037152
    "AO2" = 76.09756097560975
    "AP2" = -16.656346749226
    "AQ2" = -35.86206896551725
}

$row3Data = @{
    "A3" = 'Kazakhstan'
    "B3" = 'Kazakhstan Stock Exchange Joint-Stock Company (KAS:KASE)'
    "C3" = 'Financial Svcs. (Non-bank & Insurance)'
    "D3" = 0.197
    "E3" = 0.303
    "G3" = 0.598348623853211
    "H3" = 0.5926605504587156
    "I3" = 0.5724770642201835
    "J3" = 0.5199470768693933
    "K3" = 5.73
    "L3" = 0.5256880733944954
    "M3" = -0.0
    "N3" = -0.0
    "O3" = -0.0
    "P3" = -0.0
    "Q3" = -0.0
    "R3" = -0.0
    "S3" = 0
    "U3" = 0
    "V3" = 0
    "X3" = 0.02055145315271485
    "AB3" = 0.02055145315271485
    "AD3" = 0
    "AE3" = 0
    "AF3" = 0
    "AG3" = 0
    "AH3" = 0
    "AJ3" = 0
    "AL3" = 0.082
    "AM3" = -0.174
    "AN3" = 0
    "AO3" = 76.09756097560975
    "AP3" = 0
    "AQ3" = -35.86206896551725
}

$row4Data = @{
    "A4" = 'Kazakhstan'
    "B4" = 'Kaspi.kz Joint Stock Company (LSE:KSPI)'
    "C4" = 'Financial Svcs. (Non-bank & Insurance)'
    "D4" = 0.403
    "E4" = 0.535
    "F4" = 0.28
    "G4" = 0
    "H4" = 0
    "I4" = 0
    "J4" = 0
    "K4" = 564.9
    "L4" = 0.5595285261489699
    "M4" = 491.9
    "N4" = 0.03822037124808664
    "O4" = 0.8707735882457072
    "P4" = 491.9
    "Q4" = 0.03822037124808664
    "R4" = 0.8707735882457072
    "S4" = 0
    "T4" = 0
    "U4" = 597.3
    "V4" = 0.0464098958050054
    "W4" = 0.8295154185022026
    "X4" = 0.02089394227767173
    "Y4" = 0.8086214762245308
    "Z4" = 1.486454652532391
    "AA4" = 0
    "AB4" = 0.02196932309313021
    "AC4" = -0.02196932309313021
    "AD4" = 489.7
    "AE4" = 0
    "AF4" = 489.7
    "AG4" = -107.6
    "AH4" = 0.03665474034042426
    "AI4" = 0.4115471888393982
    "AJ4" = -0.008430950048971594
    "AK4" = -0.1815727303408707
    "AL4" = 0
    "AM4" = 0
}

foreach ($addr in $row2Data.Keys) {
    $ws.Range($addr).Value = $row2Data[$addr]
}

# B2 holds a numeric-looking label ("2") that must stay text, not become a
# number. Write it through a Text-formatted round-trip, then restore the
# default "Normal" style so no stray formatting is left behind on the cell.
$b2 = $ws.Range("B2")
$b2.NumberFormat = "@"
$b2.Value = "2"
$b2.Style = "Normal"

# T3 ("buybacks_cash_returned") is removed outright in the new data shape.
$ws.Range("T3").ClearContents()

foreach ($addr in $row3Data.Keys) {
    $ws.Range($addr).Value = $row3Data[$addr]
}

foreach ($addr in $row4Data.Keys) {
    $ws.Range($addr).Value = $row4Data[$addr]
}

"Updated rows 2-4 of the Kazakhstan financial services non-bank & insurance dataset."
